# Fruta / hortaliza, semanal
# Insert a new weekly record at row 17 (shifting existing rows 17-68 down to
# 18-69) for "Comercializadora del Agro de Limarí - Arveja Verde".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 17..68 down to 18..69, carrying formatting (date style) along.
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the new weekly record.
$ws.Cells.Item(17,1).Value = 2
$ws.Cells.Item(17,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(17,3).Value = "Coquimbo"
$ws.Cells.Item(17,4).Value = "2022-08-18"
$ws.Cells.Item(17,5).Value = 4
$ws.Cells.Item(17,6).Value = 100112022
$ws.Cells.Item(17,7).Value = "Arveja Verde"
$ws.Cells.Item(17,8).Value = "Perfection"
$ws.Cells.Item(17,9).Value = "Primera"
$ws.Cells.Item(17,10).Value = 500
$ws.Cells.Item(17,11).Value = 27000
$ws.Cells.Item(17,12).Value = 29000
$ws.Cells.Item(17,13).Value = 28000
$ws.Cells.Item(17,14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(17,15).Value = "Provincia de Limarí"
$ws.Cells.Item(17,16).Value = 1120
$ws.Cells.Item(17,17).Value = 25
$ws.Cells.Item(17,18).Value = "Hortaliza"
